$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.913.08"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "3.384.55"
$ws.Range("E3").Value = "  -2.65%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.81"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.07"
$ws.Range("E6").Value = "  +4.08%  "
$ws.Range("E7").Value = "  +4.79%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "3.385.17"
$ws.Range("E9").Value = "  -2.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.130"
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.92"
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.411"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "3.975.87"
$ws.Range("E13").Value = "  -2.60%  "
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.95"
$ws.Range("E15").Value = "  -2.96%  "
$ws.Range("D16").Value = "66.096.03"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000171"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "3.400.66"
$ws.Range("E18").Value = "  -2.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.84"
$ws.Range("E19").Value = "  -1.34%  "
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "365.43"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.50"
$ws.Range("E22").Value = "  -3.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.76"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.996"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.528"
$ws.Range("E25").Value = "  -2.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000123"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.73"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.09"
$ws.Range("E32").Value = "  -4.03%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.97"
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("E35").Value = "  -2.97%  "
$ws.Range("E36").Value = "  -2.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.59"
$ws.Range("E37").Value = "  +1.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.854"
$ws.Range("E38").Value = "  -3.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.03"
$ws.Range("E39").Value = "  -8.22%  "
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.58"
$ws.Range("E41").Value = "  +2.03%  "
$ws.Range("D42").Value = "2.672.59"
$ws.Range("E42").Value = "  -4.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.32"
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.24"
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0674"
$ws.Range("E45").Value = "  -1.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "335.45"
$ws.Range("E46").Value = "  +9.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.64"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.40"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0283"
$ws.Range("E49").Value = "  -1.37%  "
$ws.Range("E50").Value = "  +2.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "31.33"
$ws.Range("E51").Value = "  +3.50%  "
